# Append new "Print all pairs of anagrams" entry iteration to the score log.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Mark the "GFG" source for the iterations that didn't have it yet ---
$ws.Range("B100").Value = "GFG"
$ws.Range("B102").Value = "GFG"
$ws.Range("B104").Value = "GFG"
$ws.Range("B105").Value = "GFG"
$ws.Range("B106").Value = "GFG"

# --- Row 107 gains the end time in column E ---
$ws.Range("E107").Value = "5:23 - x"

# --- Row 106 gains a "done" note in column D ---
$ws.Range("D106").Value = "done 1 error, sorting strings"

# --- Row 108: tea break ---
$ws.Range("D108").Value = "tea break"
$ws.Range("E108").Value = "x - 6:18"

# --- Row 109: new topic, with hyperlink in column A ---
$ws.Range("C109").Value = "Check if characters of a given string can be rearranged to form a palindrome"
$ws.Hyperlinks.Add($ws.Range("A109"), "https://www.geeksforgeeks.org/check-characters-given-string-can-rearranged-form-palindrome/")
$ws.Range("A109").Style = "Hyperlink"
$ws.Range("E109").Value = "6:25 - x"

# --- Update selection to reflect where the user ended up editing ---
$ws.Range("E110").Select()
